$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B45").Value = "Pulkit"
$ws.Range("C45").Value = "PGI20AD010"
$ws.Range("B77").Value = "Raghav Somani"
$ws.Range("C77").Value = "PGI20CD011"
$ws.Range("B91").Value = "Subha"
$ws.Range("C91").Value = "PIET20AD054"
$ws.Range("B93").Value = "Yaduttam Pareek"
$ws.Range("C93").Value = "PGI20CS086"
$ws.Range("B96").Value = "Prem Kumar"
$ws.Range("C96").Value = "PGI20AD009"
$ws.Range("B97").Value = "Vatsal Jangid"
$ws.Range("C97").Value = "PIET20AD059"
$ws.Range("B98").Value = "Priyanshi Goyal"
$ws.Range("C98").Value = "PIET20AD039"
$ws.Range("B107").Value = "Joshi Dhiraj"
$ws.Range("C107").Value = "PIET20AD022"
$ws.Range("B108").Value = "Aastha Kanwar"
$ws.Range("C108").Value = "PGI20AD001"
$ws.Range("B113").Value = "Sourav Sharma"
$ws.Range("C113").Value = "PIET20AD053"
$ws.Range("B114").Value = "Rohit Singh Tanwar"
$ws.Range("C114").Value = "PIET20AD045"
$ws.Range("B117").Value = "Shivam Vashisht"
$ws.Range("C117").Value = "PGI20CS070"
$ws.Range("B118").Value = "Vishal Sharma"
$ws.Range("C118").Value = "PGI20CS085"
$ws.Range("B119").Value = "Gaurav Singh"
$ws.Range("C119").Value = "PGI20CS024"
$ws.Range("B123").Value = "Rajenra Choudhary"
$ws.Range("C123").Value = "PIET20AD043"
$ws.Range("B127").Value = "Rahul Lakhotiya"
$ws.Range("C127").Value = "PIET20AD041"
$ws.Range("B131").Value = "Ojasvi Sharma"
$ws.Range("C131").Value = "PIET20AD033"
$ws.Range("B132").Value = "Yash Nagal"
$ws.Range("C132").Value = "PGI20AD015"
$ws.Range("B136").Value = "Prakhar Jain"
$ws.Range("C136").Value = "PIET20AD035"
$ws.Range("B139").Value = "Tushar Suman"
$ws.Range("C139").Value = "PGI20CS082"
$ws.Range("B142").Value = "Rahul Luhar"
$ws.Range("C142").Value = "PGI20CS062"
$ws.Range("B143").Value = "Rahul Borana"
$ws.Range("C143").Value = "PGI20CS061"
$ws.Range("B150").Value = "Rahul Sharma"
$ws.Range("C150").Value = "PGI20CS063"
$ws.Range("B154").Value = "Manav Kumar"
$ws.Range("C154").Value = "PIET20AD029"
$ws.Range("B155").Value = "Yash Kumar Bhatia"
$ws.Range("C155").Value = "PIET20AD061"
$ws.Range("B158").Value = "Shubham Singh"
$ws.Range("C158").Value = "PIET20AD052"
$ws.Range("B160").Value = "Sharad Kumar"
$ws.Range("C160").Value = "PGI20CS069"
$ws.Range("B161").Value = "Ms. Shagun"
$ws.Range("C161").Value = "PIET20CS030"
$ws.Range("B167").Value = "Gargi"
$ws.Range("C167").Value = "PIER20CS208"
$ws.Range("B168").Value = "Raman Tank"
$ws.Range("C168").Value = "PGI20CS065"
$ws.Range("B170").Value = "Sushil Kumar"
$ws.Range("C170").Value = "PIET20AD057"
$ws.Range("B173").Value = "Ritik Sahu"
$ws.Range("C173").Value = "PIET20AD044"
$ws.Range("B175").Value = "Pragati Porwa"
$ws.Range("C175").Value = "PIET20AD034"
$ws.Range("B179").Value = "Siddhi Harsh"
$ws.Range("C179").Value = "PGI20CS077"
$ws.Range("B180").Value = "Chetan Sharma"
$ws.Range("C180").Value = "PIET20AD010"
$ws.Range("B181").Value = "Ankit Khemani"
$ws.Range("C181").Value = "PIET20AD005"
$ws.Range("B185").Value = "Shubham "
$ws.Range("C185").Value = "PGI20CS002"
$ws.Range("B187").Value = "Sumit Tripathi"
$ws.Range("C187").Value = "PGI20CS080"

$ws.Range("C1").Select()
